$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Copy($ws.Range("Z2"))
Write-Host "Z2 value before override:" $ws.Range("Z2").Value2
$ws.Range("Z2").Value2 = "NewText"
Write-Host "Z2 value after override:" $ws.Range("Z2").Value2
